# This script migrates the workbook from the old (AUD / 9-30-2019 / MSAUM)
# template onto the new (USD / 3-31-2021 / MSMUA) input template, per the
# commit message: 'new input template process such as new tab names'.

$wb = $excel.ActiveWorkbook

# ===================== Sector sheet =====================
$ws = $wb.Worksheets.Item("Sector")

# The old '[Unassigned]' sector row (row 13) no longer exists in the new template.
$ws.Rows.Item(13).Delete()

# Refresh AsOfDate / StrategyCode across the remaining data rows (2-12).
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 1).Value = "3/31/2021"
    $ws.Cells.Item($r, 2).Value = "MSMUA"
}

# Updated SimPortWeight (D) / IndexWeight (E) per sector.
$ws.Cells.Item(2, 4).Value = 0.0368403787726158
$ws.Cells.Item(2, 5).Value = 0.092436059678063
$ws.Cells.Item(3, 4).Value = 0.09587812042528819
$ws.Cells.Item(3, 5).Value = 0.0660985814489681
$ws.Cells.Item(4, 4).Value = 0.156846905704609
$ws.Cells.Item(4, 5).Value = 0.110359306949875
$ws.Cells.Item(5, 4).Value = 0
$ws.Cells.Item(5, 5).Value = 0.00251022935405877
$ws.Cells.Item(6, 4).Value = 0.0987685566128636
$ws.Cells.Item(6, 5).Value = 0.08954457073348471
$ws.Cells.Item(7, 4).Value = 0.147171183234493
$ws.Cells.Item(7, 5).Value = 0.188164239643132
$ws.Cells.Item(8, 4).Value = 0.0730501749977322
$ws.Cells.Item(8, 5).Value = 0.0776590703465664
$ws.Cells.Item(9, 4).Value = 0.292843039869761
$ws.Cells.Item(9, 5).Value = 0.235251101220946
$ws.Cells.Item(10, 4).Value = 0.0257900336337672
$ws.Cells.Item(10, 5).Value = 0.0280231271238259
$ws.Cells.Item(11, 4).Value = 0.0477916390022694
$ws.Cells.Item(11, 5).Value = 0.0344293370513692
$ws.Cells.Item(12, 4).Value = 0.0250199677466008
$ws.Cells.Item(12, 5).Value = 0.07552437644971109

# ===================== MCap sheet =====================
$ws = $wb.Worksheets.Item("MCap")

# The new template only carries 4 market-cap buckets (drops '< $2B').
$ws.Rows.Item(6).Delete()

# Refresh AsOfDate / StrategyCode / Currency and the bucket labels/values.
$ws.Cells.Item(2, 1).Value = "3/31/2021"
$ws.Cells.Item(2, 2).Value = "MSMUA"
$ws.Cells.Item(2, 3).Value = "USD"
$ws.Cells.Item(2, 4).Value = 1
$ws.Cells.Item(2, 5).Value = "> `$100B"
$ws.Cells.Item(2, 6).Value = 0.38758510319214
$ws.Cells.Item(2, 7).Value = 0.370893687603703
$ws.Cells.Item(3, 1).Value = "3/31/2021"
$ws.Cells.Item(3, 2).Value = "MSMUA"
$ws.Cells.Item(3, 3).Value = "USD"
$ws.Cells.Item(3, 4).Value = 2
$ws.Cells.Item(3, 5).Value = "`$25B - `$100B"
$ws.Cells.Item(3, 6).Value = 0.387982644642594
$ws.Cells.Item(3, 7).Value = 0.454372528558921
$ws.Cells.Item(4, 1).Value = "3/31/2021"
$ws.Cells.Item(4, 2).Value = "MSMUA"
$ws.Cells.Item(4, 3).Value = "USD"
$ws.Cells.Item(4, 4).Value = 3
$ws.Cells.Item(4, 5).Value = "`$15B - `$25B"
$ws.Cells.Item(4, 6).Value = 0.163760518702805
$ws.Cells.Item(4, 7).Value = 0.104664363648447
$ws.Cells.Item(5, 1).Value = "3/31/2021"
$ws.Cells.Item(5, 2).Value = "MSMUA"
$ws.Cells.Item(5, 3).Value = "USD"
$ws.Cells.Item(5, 4).Value = 4
$ws.Cells.Item(5, 5).Value = "`$2B - `$15B"
$ws.Cells.Item(5, 6).Value = 0.06067173346246151
$ws.Cells.Item(5, 7).Value = 0.0700694201889292

# ===================== Chars sheet =====================
$ws = $wb.Worksheets.Item("Chars")

# Refresh AsOfDate / StrategyCode / Currency and the portfolio + index stats.
$ws.Cells.Item(2, 1).Value = "3/31/2021"
$ws.Cells.Item(2, 2).Value = "MSMUA"
$ws.Cells.Item(2, 3).Value = "USD"
$ws.Cells.Item(2, 4).Value = "Number of Securities"
$ws.Cells.Item(2, 5).Value = 79
$ws.Cells.Item(2, 6).Value = 185
$ws.Cells.Item(3, 1).Value = "3/31/2021"
$ws.Cells.Item(3, 2).Value = "MSMUA"
$ws.Cells.Item(3, 3).Value = "USD"
$ws.Cells.Item(3, 4).Value = "Price/Earnings Ratio (LTM)"
$ws.Cells.Item(3, 5).Value = "30.79"
$ws.Cells.Item(3, 6).Value = "31.49"
$ws.Cells.Item(4, 1).Value = "3/31/2021"
$ws.Cells.Item(4, 2).Value = "MSMUA"
$ws.Cells.Item(4, 3).Value = "USD"
$ws.Cells.Item(4, 4).Value = "Dividend Yield (Current)"
$ws.Cells.Item(4, 5).Value = "1.87%"
$ws.Cells.Item(4, 6).Value = "2.03%"
$ws.Cells.Item(5, 1).Value = "3/31/2021"
$ws.Cells.Item(5, 2).Value = "MSMUA"
$ws.Cells.Item(5, 3).Value = "USD"
$ws.Cells.Item(5, 4).Value = "EPS Growth (5 Yr. Historical)"
$ws.Cells.Item(5, 5).Value = "12.17%"
$ws.Cells.Item(5, 6).Value = "11.92%"
$ws.Cells.Item(6, 1).Value = "3/31/2021"
$ws.Cells.Item(6, 2).Value = "MSMUA"
$ws.Cells.Item(6, 3).Value = "USD"
$ws.Cells.Item(6, 4).Value = "Price/Book Ratio"
$ws.Cells.Item(6, 5).Value = "6.59"
$ws.Cells.Item(6, 6).Value = "6.47"
$ws.Cells.Item(7, 1).Value = "3/31/2021"
$ws.Cells.Item(7, 2).Value = "MSMUA"
$ws.Cells.Item(7, 3).Value = "USD"
$ws.Cells.Item(7, 4).Value = "Weighted Average Market Cap"
$ws.Cells.Item(7, 5).Value = "`$101.5 B"
$ws.Cells.Item(7, 6).Value = "`$146.1 B"
$ws.Cells.Item(8, 1).Value = "3/31/2021"
$ws.Cells.Item(8, 2).Value = "MSMUA"
$ws.Cells.Item(8, 3).Value = "USD"
$ws.Cells.Item(8, 4).Value = "Weighted Median Market Cap"
$ws.Cells.Item(8, 5).Value = "`$62.3 B"
$ws.Cells.Item(8, 6).Value = "`$62.0 B"

# ===================== TB sheet (top holdings) =====================
$ws = $wb.Worksheets.Item("TB")

# Entire top-10 holdings list is replaced for the new (US-equity) template.
$ws.Cells.Item(2, 1).Value = "3/31/2021"
$ws.Cells.Item(2, 2).Value = "MSMUA"
$ws.Cells.Item(2, 3).Value = "Accenture Plc Class A"
$ws.Cells.Item(2, 4).Value = 0.0581071909668409
$ws.Cells.Item(2, 5).Value = 0.3433366888372114
$ws.Cells.Item(3, 1).Value = "3/31/2021"
$ws.Cells.Item(3, 2).Value = "MSMUA"
$ws.Cells.Item(3, 3).Value = "Progressive Corporation"
$ws.Cells.Item(3, 4).Value = 0.0468109429178378
$ws.Cells.Item(3, 5).Value = 0.3433366888372114
$ws.Cells.Item(4, 1).Value = "3/31/2021"
$ws.Cells.Item(4, 2).Value = "MSMUA"
$ws.Cells.Item(4, 3).Value = "Bristol-Myers Squibb Company"
$ws.Cells.Item(4, 4).Value = 0.04056710513628339
$ws.Cells.Item(4, 5).Value = 0.3433366888372114
$ws.Cells.Item(5, 1).Value = "3/31/2021"
$ws.Cells.Item(5, 2).Value = "MSMUA"
$ws.Cells.Item(5, 3).Value = "Crown Castle International Corp"
$ws.Cells.Item(5, 4).Value = 0.0374801068532375
$ws.Cells.Item(5, 5).Value = 0.3433366888372114
$ws.Cells.Item(6, 1).Value = "3/31/2021"
$ws.Cells.Item(6, 2).Value = "MSMUA"
$ws.Cells.Item(6, 3).Value = "International Business Machines Corporation"
$ws.Cells.Item(6, 4).Value = 0.0314563481697103
$ws.Cells.Item(6, 5).Value = 0.3433366888372114
$ws.Cells.Item(7, 1).Value = "3/31/2021"
$ws.Cells.Item(7, 2).Value = "MSMUA"
$ws.Cells.Item(7, 3).Value = "Johnson & Johnson"
$ws.Cells.Item(7, 4).Value = 0.0268942334700008
$ws.Cells.Item(7, 5).Value = 0.3433366888372114
$ws.Cells.Item(8, 1).Value = "3/31/2021"
$ws.Cells.Item(8, 2).Value = "MSMUA"
$ws.Cells.Item(8, 3).Value = "McCormick & Company, Incorporated"
$ws.Cells.Item(8, 4).Value = 0.0264738775025792
$ws.Cells.Item(8, 5).Value = 0.3433366888372114
$ws.Cells.Item(9, 1).Value = "3/31/2021"
$ws.Cells.Item(9, 2).Value = "MSMUA"
$ws.Cells.Item(9, 3).Value = "Target Corporation"
$ws.Cells.Item(9, 4).Value = 0.0261026620817259
$ws.Cells.Item(9, 5).Value = 0.3433366888372114
$ws.Cells.Item(10, 1).Value = "3/31/2021"
$ws.Cells.Item(10, 2).Value = "MSMUA"
$ws.Cells.Item(10, 3).Value = "Citrix Systems, Inc."
$ws.Cells.Item(10, 4).Value = 0.0247964746261889
$ws.Cells.Item(10, 5).Value = 0.3433366888372114
$ws.Cells.Item(11, 1).Value = "3/31/2021"
$ws.Cells.Item(11, 2).Value = "MSMUA"
$ws.Cells.Item(11, 3).Value = "Campbell Soup Company"
$ws.Cells.Item(11, 4).Value = 0.0246477471128067
$ws.Cells.Item(11, 5).Value = 0.3433366888372114
